$d = $word.ActiveDocument

$replacements = @(
    @("2025-11-20 Thursday", "2025-11-21 Friday"),
    @("29×16=464", "88×12=1056"),
    @("70×76=5320", "53×28=1484"),
    @("61×98=5978", "24×23=552"),
    @("28×56=1568", "81×74=5994"),
    @("59×77=4543", "11×82=902"),
    @("14×70=980", "70×41=2870"),
    @("59×99=5841", "39×12=468"),
    @("50×24=1200", "79×28=2212"),
    @("83×90=7470", "21×43=903"),
    @("55×81=4455", "31×17=527"),
    @("16×35=560", "72×14=1008"),
    @("90×84=7560", "25×17=425"),
    @("83×78=6474", "16×51=816"),
    @("42×51=2142", "94×18=1692"),
    @("62×28=1736", "45×43=1935"),
    @("76×51=3876", "53×23=1219"),
    @("24×92=2208", "20×98=1960"),
    @("25×77=1925", "43×96=4128"),
    @("78×95=7410", "87×61=5307"),
    @("75×78=5850", "54×71=3834"),
    @("54×73=3942", "35×15=525"),
    @("67×35=2345", "30×95=2850"),
    @("22×38=836", "79×67=5293"),
    @("72×22=1584", "68×95=6460"),
    @("65×76=4940", "51×15=765")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
